$wb = $excel.ActiveWorkbook

# --- Sheet "survey": insert a new column D (inputAttributes.data-type) ---
$survey = $wb.Worksheets.Item("survey")

# Insert a new, blank column at D. This shifts the existing
# appearance/condition/name/label data (old D/E/F) one column to the
# right (new E/F/G), preserving their original widths untouched.
$survey.Columns.Item(4).Insert()

# New column header + two new appearance/inputAttributes.data-type values
# for the new rows below.
$survey.Range("D1").Value = "inputAttributes.data-type"

# New row 16: a horizontal select_one yes_no example
$survey.Range("B16").Value = "select_one yes_no"
$survey.Range("D16").Value = "horizontal"
$survey.Range("F16").Value = "h_select"
$survey.Range("G16").Value = "Horizontal select example."

# New row 17: a select backed by a content-provider query
$survey.Range("B17").Value = "select_one content_provider_test"
$survey.Range("F17").Value = "cp_test"
$survey.Range("G17").Value = "This demos a content provider query."

# Column width adjustments to match the authored layout.
$survey.Columns.Item(4).ColumnWidth = 20.71
$survey.Columns.Item(7).ColumnWidth = 46.0

# --- Sheet "queries": rename the odk_values query, add a literal result ---
$queries = $wb.Worksheets.Item("queries")
$queries.Range("A5").Value = "content_provider_test"
$queries.Range("C5").Value = '[{ name: "test", label : JSON.stringify(context) }]'
